$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update AgTests (F) and AgPosit (G) values for rows 272-433
$ws.Cells.Item(272, 6).Value = 30512
$ws.Cells.Item(274, 6).Value = 28148
$ws.Cells.Item(275, 6).Value = 30369
$ws.Cells.Item(276, 6).Value = 11433
$ws.Cells.Item(277, 6).Value = 3384
$ws.Cells.Item(278, 6).Value = 30544
$ws.Cells.Item(280, 6).Value = 34810
$ws.Cells.Item(281, 6).Value = 46079
$ws.Cells.Item(282, 6).Value = 46295
$ws.Cells.Item(284, 6).Value = 1197
$ws.Cells.Item(285, 6).Value = 42102
$ws.Cells.Item(288, 6).Value = 59270
$ws.Cells.Item(289, 6).Value = 63023
$ws.Cells.Item(290, 6).Value = 17613
$ws.Cells.Item(291, 6).Value = 15144
$ws.Cells.Item(294, 6).Value = 93926
$ws.Cells.Item(296, 6).Value = 1851
$ws.Cells.Item(297, 6).Value = 2399
$ws.Cells.Item(298, 6).Value = 3241
$ws.Cells.Item(299, 6).Value = 65690
$ws.Cells.Item(301, 6).Value = 72220
$ws.Cells.Item(302, 6).Value = 78587
$ws.Cells.Item(304, 6).Value = 6105
$ws.Cells.Item(305, 6).Value = 3368
$ws.Cells.Item(306, 6).Value = 75359
$ws.Cells.Item(308, 6).Value = 15562
$ws.Cells.Item(309, 6).Value = 77902
$ws.Cells.Item(310, 6).Value = 79247
$ws.Cells.Item(311, 6).Value = 61499
$ws.Cells.Item(312, 6).Value = 28177
$ws.Cells.Item(313, 6).Value = 75814
$ws.Cells.Item(315, 6).Value = 56390
$ws.Cells.Item(316, 6).Value = 50726
$ws.Cells.Item(317, 6).Value = 63745
$ws.Cells.Item(318, 6).Value = 49342
$ws.Cells.Item(319, 6).Value = 41371
$ws.Cells.Item(320, 6).Value = 73448
$ws.Cells.Item(320, 7).Value = 3359
$ws.Cells.Item(321, 6).Value = 89816
$ws.Cells.Item(321, 7).Value = 2668
$ws.Cells.Item(322, 6).Value = 109384
$ws.Cells.Item(323, 6).Value = 217573
$ws.Cells.Item(324, 6).Value = 250040
$ws.Cells.Item(324, 7).Value = 2859
$ws.Cells.Item(325, 6).Value = 775341
$ws.Cells.Item(325, 7).Value = 6517
$ws.Cells.Item(326, 6).Value = 418135
$ws.Cells.Item(327, 6).Value = 223825
$ws.Cells.Item(328, 6).Value = 181106
$ws.Cells.Item(329, 6).Value = 73375
$ws.Cells.Item(330, 6).Value = 71581
$ws.Cells.Item(330, 7).Value = 2079
$ws.Cells.Item(331, 6).Value = 153957
$ws.Cells.Item(331, 7).Value = 2708
$ws.Cells.Item(333, 6).Value = 254959
$ws.Cells.Item(334, 6).Value = 192978
$ws.Cells.Item(335, 6).Value = 150233
$ws.Cells.Item(335, 7).Value = 3748
$ws.Cells.Item(336, 6).Value = 81688
$ws.Cells.Item(336, 7).Value = 2562
$ws.Cells.Item(337, 6).Value = 104005
$ws.Cells.Item(337, 7).Value = 2891
$ws.Cells.Item(338, 6).Value = 221364
$ws.Cells.Item(338, 7).Value = 3047
$ws.Cells.Item(339, 6).Value = 662660
$ws.Cells.Item(340, 6).Value = 387179
$ws.Cells.Item(340, 7).Value = 3308
$ws.Cells.Item(341, 6).Value = 283345
$ws.Cells.Item(341, 7).Value = 3609
$ws.Cells.Item(343, 6).Value = 133371
$ws.Cells.Item(344, 6).Value = 135568
$ws.Cells.Item(345, 6).Value = 292202
$ws.Cells.Item(346, 6).Value = 674718
$ws.Cells.Item(347, 6).Value = 346587
$ws.Cells.Item(347, 7).Value = 2919
$ws.Cells.Item(348, 6).Value = 232778
$ws.Cells.Item(350, 6).Value = 127010
$ws.Cells.Item(351, 6).Value = 150482
$ws.Cells.Item(352, 6).Value = 307360
$ws.Cells.Item(353, 6).Value = 723556
$ws.Cells.Item(354, 6).Value = 316744
$ws.Cells.Item(354, 7).Value = 2883
$ws.Cells.Item(355, 6).Value = 221948
$ws.Cells.Item(357, 6).Value = 138203
$ws.Cells.Item(358, 6).Value = 158772
$ws.Cells.Item(359, 6).Value = 321158
$ws.Cells.Item(360, 6).Value = 749715
$ws.Cells.Item(361, 6).Value = 332874
$ws.Cells.Item(362, 6).Value = 228546
$ws.Cells.Item(364, 6).Value = 168406
$ws.Cells.Item(365, 6).Value = 184600
$ws.Cells.Item(366, 6).Value = 339402
$ws.Cells.Item(367, 6).Value = 767082
$ws.Cells.Item(367, 7).Value = 3924
$ws.Cells.Item(368, 6).Value = 346192
$ws.Cells.Item(368, 7).Value = 2300
$ws.Cells.Item(369, 6).Value = 234701
$ws.Cells.Item(370, 6).Value = 180882
$ws.Cells.Item(371, 6).Value = 160007
$ws.Cells.Item(372, 6).Value = 178380
$ws.Cells.Item(373, 6).Value = 350030
$ws.Cells.Item(374, 6).Value = 773696
$ws.Cells.Item(375, 6).Value = 351270
$ws.Cells.Item(375, 7).Value = 1857
$ws.Cells.Item(376, 6).Value = 221397
$ws.Cells.Item(377, 6).Value = 176549
$ws.Cells.Item(378, 6).Value = 157269
$ws.Cells.Item(379, 6).Value = 179528
$ws.Cells.Item(380, 6).Value = 344425
$ws.Cells.Item(381, 6).Value = 746391
$ws.Cells.Item(381, 7).Value = 2696
$ws.Cells.Item(383, 6).Value = 220795
$ws.Cells.Item(384, 6).Value = 171588
$ws.Cells.Item(385, 6).Value = 150753
$ws.Cells.Item(386, 6).Value = 182756
$ws.Cells.Item(387, 6).Value = 351555
$ws.Cells.Item(388, 6).Value = 730436
$ws.Cells.Item(388, 7).Value = 2204
$ws.Cells.Item(390, 6).Value = 219782
$ws.Cells.Item(391, 6).Value = 177253
$ws.Cells.Item(392, 6).Value = 220876
$ws.Cells.Item(393, 6).Value = 308756
$ws.Cells.Item(393, 7).Value = 1242
$ws.Cells.Item(394, 6).Value = 166545
$ws.Cells.Item(395, 6).Value = 751888
$ws.Cells.Item(395, 7).Value = 1964
$ws.Cells.Item(396, 6).Value = 164905
$ws.Cells.Item(397, 6).Value = 108107
$ws.Cells.Item(398, 6).Value = 299014
$ws.Cells.Item(398, 7).Value = 1472
$ws.Cells.Item(399, 6).Value = 201578
$ws.Cells.Item(399, 7).Value = 967
$ws.Cells.Item(400, 6).Value = 151150
$ws.Cells.Item(400, 7).Value = 771
$ws.Cells.Item(401, 6).Value = 273695
$ws.Cells.Item(402, 6).Value = 718231
$ws.Cells.Item(403, 6).Value = 352320
$ws.Cells.Item(403, 7).Value = 734
$ws.Cells.Item(404, 6).Value = 225168
$ws.Cells.Item(404, 7).Value = 913
$ws.Cells.Item(405, 6).Value = 174638
$ws.Cells.Item(406, 6).Value = 170666
$ws.Cells.Item(407, 6).Value = 158355
$ws.Cells.Item(408, 6).Value = 304524
$ws.Cells.Item(409, 6).Value = 704667
$ws.Cells.Item(410, 6).Value = 363888
$ws.Cells.Item(410, 7).Value = 630
$ws.Cells.Item(411, 6).Value = 225036
$ws.Cells.Item(412, 6).Value = 175992
$ws.Cells.Item(412, 7).Value = 646
$ws.Cells.Item(413, 6).Value = 149196
$ws.Cells.Item(414, 6).Value = 147155
$ws.Cells.Item(415, 6).Value = 304906
$ws.Cells.Item(415, 7).Value = 693
$ws.Cells.Item(416, 6).Value = 660384
$ws.Cells.Item(417, 6).Value = 332988
$ws.Cells.Item(418, 6).Value = 200837
$ws.Cells.Item(419, 6).Value = 147637
$ws.Cells.Item(419, 7).Value = 504
$ws.Cells.Item(420, 6).Value = 137008
$ws.Cells.Item(420, 7).Value = 492
$ws.Cells.Item(421, 6).Value = 150994
$ws.Cells.Item(421, 7).Value = 526
$ws.Cells.Item(422, 6).Value = 294542
$ws.Cells.Item(422, 7).Value = 642
$ws.Cells.Item(423, 6).Value = 432400
$ws.Cells.Item(423, 7).Value = 629
$ws.Cells.Item(425, 6).Value = 137628
$ws.Cells.Item(426, 6).Value = 105690
$ws.Cells.Item(426, 7).Value = 385
$ws.Cells.Item(427, 6).Value = 89868
$ws.Cells.Item(427, 7).Value = 363
$ws.Cells.Item(428, 6).Value = 100827
$ws.Cells.Item(428, 7).Value = 385
$ws.Cells.Item(429, 6).Value = 174085
$ws.Cells.Item(429, 7).Value = 450
$ws.Cells.Item(430, 6).Value = 170990
$ws.Cells.Item(431, 6).Value = 163793
$ws.Cells.Item(431, 7).Value = 392
$ws.Cells.Item(432, 6).Value = 121927
$ws.Cells.Item(432, 7).Value = 424
$ws.Cells.Item(433, 6).Value = 84555
$ws.Cells.Item(433, 7).Value = 259

# Append new row 434 with data for 2021-05-12
$ws.Cells.Item(434, 1).Value = 44328
$ws.Cells.Item(434, 2).Value = 386868
$ws.Cells.Item(434, 3).Value = 4968
$ws.Cells.Item(434, 4).Value = 328
$ws.Cells.Item(434, 5).Value = 12135
$ws.Cells.Item(434, 6).Value = 64512
$ws.Cells.Item(434, 7).Value = 219

# Ensure A434 uses the same date number format as the rest of column A
$ws.Cells.Item(434, 1).NumberFormat = "yyyy-mm-dd"
